# Added tests for Gradebook class
# - Updates row 7 (Amit/Freifeld, perfect score) on "grades"
# - Appends rows 8-19 with additional test students covering each
#   weight individually plus a few blended/partial-credit scenarios
# - Recomputes the SUMPRODUCT/VLOOKUP formulas for every new row
# - Moves the active selection to C5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grades")

# --- Row 7: repurposed as "Amit" / "Freifeld", full marks -------------
$ws.Range("A7").Value = "Amit"
$ws.Range("B7").Value = "Freifeld"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("H7").Formula = '=SUMPRODUCT(C7:F7,finalscores!A$1:D$1)'
$ws.Range("I7").Formula = '=VLOOKUP(H7,gradeboundaries,2)'

# --- New rows 8-19: synthetic test fixtures ----------------------------
$rows = @(
    @{ r = 8;  a = "Clark ";  b = "Shesh";   c = 1;    d = 1;    e = 1;    f = 1 },
    @{ r = 9;  a = "Bob";     b = "Builder";  c = 0;    d = 0;    e = 0;    f = 0 },
    @{ r = 10; a = "Fname1";  b = "Lname1";   c = 1;    d = 0;    e = 0;    f = 0 },
    @{ r = 11; a = "Fname2";  b = "Lname2";   c = 0;    d = 1;    e = 0;    f = 0 },
    @{ r = 12; a = "Fname3";  b = "Lname3";   c = 0;    d = 0;    e = 1;    f = 0 },
    @{ r = 13; a = "Fname4";  b = "Lname4";   c = 0;    d = 0;    e = 0;    f = 1 },
    @{ r = 14; a = "Fname5";  b = "Lname5";   c = 0.2;  d = 1;    e = 1;    f = 1 },
    @{ r = 15; a = "Fname6";  b = "Lname6";   c = 1;    d = 0.4;  e = 1;    f = 1 },
    @{ r = 16; a = "Fname7";  b = "Lname7";   c = 1;    d = 1;    e = 0.4;  f = 1 },
    @{ r = 17; a = "Fname8";  b = "Lname8";   c = 1;    d = 0.7;  e = 1;    f = 0.1 },
    @{ r = 18; a = "Fname9";  b = "Lname9";   c = 0.34; d = 1;    e = 0.4;  f = 1 },
    @{ r = 19; a = "Fname10"; b = "Lname10";  c = 1;    d = 0.34; e = 1;    f = 1 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.a
    $ws.Range("B$r").Value = $row.b
    $ws.Range("C$r").Value = $row.c
    $ws.Range("D$r").Value = $row.d
    $ws.Range("E$r").Value = $row.e
    $ws.Range("F$r").Value = $row.f
    $ws.Range("G$r").Value = ""
    $ws.Range("H$r").Formula = "=SUMPRODUCT(C${r}:F${r},finalscores!A`$1:D`$1)"
    $ws.Range("I$r").Formula = "=VLOOKUP(H$r,gradeboundaries,2)"
}

# Row 20 stays present but blank (keeps the dimension in sync with the diff)
$ws.Rows.Item(20).RowHeight = 13.8

# --- Selection moves to C5 ---------------------------------------------
$ws.Range("C5").Select()
